{"js": "// Grammar-correction edit for legendary-robot README document.\n// Each fix below locates the exact old text with body.search() (literal,\n// case-sensitive match) and rewrites it in place via insertText(..., \"Replace\").\n\nasync function replaceOnce(body, oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. Turn \"Node.js Download Page\" into prose with an inline markdown link to\n//    the Node.js download page (and keep the two-space gap before the next\n//    \"{0}{1} 2.\" line-break marker that the author used elsewhere).\nawait replaceOnce(\n  body,\n  \"{0}{1} 1. Node.js Download Page {0}{1} 2. After install, check with command line to ensure setup is correct.\",\n  \"{0}{1} 1. Node.js download page [link](https://nodejs.org/en/download/)  {0}{1} 2. After install, check with command line to ensure setup is correct.\"\n);\n\n// 2. USAGE section reference to the \"js\" folder: drop the redundant \"in repo\".\nawait replaceOnce(\n  body,\n  \"\\u201d folder in repo and install the \",\n  \"\\u201d folder and install the \"\n);\n\n// 3. Rephrase intro sentence of the USAGE section.\nawait replaceOnce(\n  body,\n  \"Follow the below steps once you've completed the installation section guidelines.\",\n  \"Once you've completed the installation section guidelines, follow the below steps to run the application.\"\n);\n\n// 4. Add the \"node index.js\" instruction and drop \"in repo\" again.\nawait replaceOnce(\n  body,\n  \"\\u201d folder in repo. See below image for guidance.\",\n  \"\\u201d folder and type `node index.js`. See below image for guidance.\"\n);\n\n// 5. Remove the duplicated \"Clone repository and open in VS Code.\" step (it\n//    is already covered in the INSTALLATION section) and renumber step 3.\nawait replaceOnce(\n  body,\n  \"{0}{1} 3. Clone repository and open in VS Code. Once you hit enter, the application will begin to\",\n  \"{0}{1} 3. Once you hit enter, the application will begin to\"\n);\n\n// 6. Grammar: \"answer the prompts\" -> \"answering prompts\".\nawait replaceOnce(\n  body,\n  \"when answer the prompts.\",\n  \"when answering prompts.\"\n);\n\n// 7. Grammar: \"used for generating\" -> \"used to generate\".\nawait replaceOnce(\n  body,\n  \"Check out the following link for the text that was used for generating this README document\",\n  \"Check out the following link for the text that was used to generate this README document\"\n);\n\n// 8. Add a \"###\" markdown heading marker in front of \"VS Code Preview\".\nawait replaceOnce(\n  body,\n  \"numbers.{0}{0} VS Code Preview \",\n  \"numbers.{0}{0} ### VS Code Preview \"\n);\n\n// 9. Rephrase the VS Code preview-feature sentence (grammar + tense fixes).\nawait replaceOnce(\n  body,\n  \"One efficiency tips when writing README files was utilizing VS Codes preview markdown feature. This feature allows you to preview how your markdown files before having to upload to see the final output.\",\n  \"One efficiency tips when writing README files is to utilize VS Codes preview markdown feature. This feature allows you to preview how your markdown file looks before having to upload to see the final output.\"\n);\n", "ps1": "# Grammar-correction edit for legendary-robot README document.\n# Each fix below locates the exact old text with Range.Find.Execute (literal,\n# case-sensitive match) and rewrites the matched Range's Text directly.\n# (Assigning Find.Execute's own ReplaceWith/Replace parameters routes through\n# AutoFormat and silently turns straight apostrophes into curly ones; setting\n# Range.Text on the found hit keeps the literal characters we ask for.)\n\nfunction Replace-Once($doc, $oldText, $newText) {\n    $range = $doc.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n    $range.Text = $newText\n}\n\n$d = $word.ActiveDocument\n\n# 1. Turn \"Node.js Download Page\" into prose with an inline markdown link to\n#    the Node.js download page (and keep the two-space gap before the next\n#    \"{0}{1} 2.\" line-break marker that the author used elsewhere).\nReplace-Once $d \"{0}{1} 1. Node.js Download Page {0}{1} 2. After install, check with command line to ensure setup is correct.\" \"{0}{1} 1. Node.js download page [link](https://nodejs.org/en/download/)  {0}{1} 2. After install, check with command line to ensure setup is correct.\"\n\n# 2. USAGE section reference to the \"js\" folder: drop the redundant \"in repo\".\nReplace-Once $d \"\u201d folder in repo and install the \" \"\u201d folder and install the \"\n\n# 3. Rephrase intro sentence of the USAGE section.\nReplace-Once $d \"Follow the below steps once you've completed the installation section guidelines.\" \"Once you've completed the installation section guidelines, follow the below steps to run the application.\"\n\n# 4. Add the \"node index.js\" instruction and drop \"in repo\" again.\nReplace-Once $d \"\u201d folder in repo. See below image for guidance.\" \"\u201d folder and type ``node index.js``. See below image for guidance.\"\n\n# 5. Remove the duplicated \"Clone repository and open in VS Code.\" step (it\n#    is already covered in the INSTALLATION section) and renumber step 3.\nReplace-Once $d \"{0}{1} 3. Clone repository and open in VS Code. Once you hit enter, the application will begin to\" \"{0}{1} 3. Once you hit enter, the application will begin to\"\n\n# 6. Grammar: \"answer the prompts\" -> \"answering prompts\".\nReplace-Once $d \"when answer the prompts.\" \"when answering prompts.\"\n\n# 7. Grammar: \"used for generating\" -> \"used to generate\".\nReplace-Once $d \"Check out the following link for the text that was used for generating this README document\" \"Check out the following link for the text that was used to generate this README document\"\n\n# 8. Add a \"###\" markdown heading marker in front of \"VS Code Preview\".\nReplace-Once $d \"numbers.{0}{0} VS Code Preview \" \"numbers.{0}{0} ### VS Code Preview \"\n\n# 9. Rephrase the VS Code preview-feature sentence (grammar + tense fixes).\nReplace-Once $d \"One efficiency tips when writing README files was utilizing VS Codes preview markdown feature. This feature allows you to preview how your markdown files before having to upload to see the final output.\" \"One efficiency tips when writing README files is to utilize VS Codes preview markdown feature. This feature allows you to preview how your markdown file looks before having to upload to see the final output.\"\n"}
